# Edits PlayerPerformance_4730.xlsx:
#  1) Insert a new "Player Info" sheet before the existing "ODI Batting" sheet
#     with player ID/NAME/BATTING_HAND/BOWL_STYLE details.
#  2) On the "ODI Batting" sheet, rename the MATCH_CARD_LINK column to
#     MATCH_CODE and replace the full scorecard URLs with the bare match
#     code that was embedded in them.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) New "Player Info" sheet
# ---------------------------------------------------------------------------
# Worksheets.Add() inserts the new sheet ahead of the currently active sheet,
# which puts it first - matching the target sheet order (Player Info, then
# ODI Batting).
$playerInfo = $wb.Worksheets.Add()
$playerInfo.Name = "Player Info"

$playerInfo.Range("A1").Value = "ID"
$playerInfo.Range("B1").Value = "NAME"
$playerInfo.Range("C1").Value = "BATTING_HAND"
$playerInfo.Range("D1").Value = "BOWL_STYLE"

$piHeader = $playerInfo.Range("A1:D1")
$piHeader.Font.Bold = $true
$piHeader.Borders.LineStyle = 1
$piHeader.HorizontalAlignment = -4108
$piHeader.VerticalAlignment = -4160

# Keep the ID as text (e.g. "4730") instead of a number.
$playerInfo.Range("A2").NumberFormat = "@"
$playerInfo.Range("A2").Value = "4730"
$playerInfo.Range("B2").Value = "Khayelihle Zondo"
$playerInfo.Range("C2").Value = "Right Handed"
$playerInfo.Range("D2").Value = "Right Arm Off Break"

# ---------------------------------------------------------------------------
# 2) "ODI Batting" sheet: MATCH_CARD_LINK -> MATCH_CODE
# ---------------------------------------------------------------------------
$odiBatting = $wb.Worksheets.Item("ODI Batting")

$odiBatting.Range("D1").Value = "MATCH_CODE"

$matchCodes = @{
    2 = "4127"
    3 = "4128"
    4 = "4135"
    5 = "4207"
    6 = "4208"
    7 = "4517"
}

foreach ($row in $matchCodes.Keys) {
    $cell = $odiBatting.Range("D$row")
    $cell.NumberFormat = "@"
    $cell.Value = $matchCodes[$row]
}
